$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "2025/12/03 02:00"
$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = "-"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = "-"
